$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "autumn" row's complementary color (G2) was updated.
$ws.Range("G2").Value = "#FB9637"

# Reflect the new selection location recorded in the saved view state.
$ws.Range("G2").Select()
